# Updated cryptos list on Fri Sep 29 17:47:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.869.21"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "1.666.50"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "215.46"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +5.09%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "20.18"
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "0.0896"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "1.901.44"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "1.653.08"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "66.03"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "26.880.69"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "232.10"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("D19").Value = "7.80"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").Value = "145.55"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "7.12"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "1.464.89"
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("D34").Value = "3.16"
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.899"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.573"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  +6.70%  "
$ws.Range("D44").Value = "65.92"
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").Value = "1.813.05"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").Value = "0.777"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "90.25"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").Value = "0.0507"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "7.59"
$ws.Range("E51").Value = "  +0.71%  "
